$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet references ("Control Box Physical Interface", "Relay Box Physical
# Interface", "Radio Box Physical Interface" -> rId1/rId2/rId3 -> index 1/2/3)
# ---------------------------------------------------------------------------
$wsControl = $wb.Worksheets.Item(1)
$wsRelay   = $wb.Worksheets.Item(2)
$wsRadio   = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# New connector: "Air Horn" added to the Control Box Physical Interface sheet
# as rows 42-43, mirroring the layout/style of the other merged interface
# blocks in column A (e.g. rows 4-5, 6-8, ...).
# ---------------------------------------------------------------------------
$wsControl.Range("A42").Value = "Air Horn"
$wsControl.Range("B42").Value = "M12A-2 female"
$wsControl.Range("C42").Value = 1
$wsControl.Range("D42").Value = "Horn"
$wsControl.Range("E42").Value = "12V"
$wsControl.Range("F42").Value = "power"
$wsControl.Range("G42").Value = "Boatse 0.2"

$wsControl.Range("B43").Value = "M12A-2 female"
$wsControl.Range("C43").Value = 1
$wsControl.Range("D43").Value = "GND"
$wsControl.Range("E43").Value = "GND"
$wsControl.Range("F43").Value = "ground"
$wsControl.Range("G43").Value = "Boatse 0.2"

# Match the formatting used by the other merged "Interface" column cells:
# centered horizontally/vertically, merged across the two connector rows.
$newBlock = $wsControl.Range("A42:A43")
$newBlock.HorizontalAlignment = -4108
$newBlock.VerticalAlignment = -4108
$newBlock.Merge()

# ---------------------------------------------------------------------------
# View/selection state.
# The active tab moves from "Radio Box Physical Interface" (sheet 3) to
# "Control Box Physical Interface" (sheet 1), and each sheet's remembered
# selection changes.
# ---------------------------------------------------------------------------
$wsRelay.Activate()
$wsRelay.Range("F17").Select()

$wsRadio.Activate()
$wsRadio.Range("B7").Select()

$wsControl.Activate()
$wsControl.Range("B6").Select()
